# precision-recall-graph.xlsx : update the cached precision values on the
# "Data" sheet (normalization of the underlying vectors), and move the
# active selection to D15 on Data (which becomes the active sheet instead
# of the Chart1 chartsheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# --- Column B ("Normal") ---
$ws.Range("B3").Value = 0.73653905116585106
$ws.Range("B4").Value = 0.734720869347669
$ws.Range("B5").Value = 0.726720869347669
$ws.Range("B6").Value = 0.70549864712544696
$ws.Range("B7").Value = 0.680215097341897
$ws.Range("B8").Value = 0.66947700210380201
$ws.Range("B9").Value = 0.58250201203219099
$ws.Range("B10").Value = 0.56217381408634604
$ws.Range("B11").Value = 0.53421701752831596
$ws.Range("B12").Value = 0.50966259624643095
$ws.Range("B13").Value = 0.50842283016455903

# --- Column C ("Stopwords") ---
$ws.Range("C3").Value = 0.73170126930228696
$ws.Range("C4").Value = 0.72988308748410602
$ws.Range("C5").Value = 0.72719336416394798
$ws.Range("C6").Value = 0.70386003083061399
$ws.Range("C7").Value = 0.68179726026784404
$ws.Range("C8").Value = 0.67100629036817805
$ws.Range("C9").Value = 0.58693483653820999
$ws.Range("C10").Value = 0.56007358177989697
$ws.Range("C11").Value = 0.53343777477644305
$ws.Range("C12").Value = 0.50660730127306797
$ws.Range("C13").Value = 0.50548031714608399

# --- Column D ("Stemming") ---
$ws.Range("D3").Value = 0.74313543198619703
$ws.Range("D4").Value = 0.74113543198619702
$ws.Range("D5").Value = 0.73446876531952998
$ws.Range("D6").Value = 0.70002432087508604
$ws.Range("D7").Value = 0.69269098754175196
$ws.Range("D8").Value = 0.68801252590676898
$ws.Range("D9").Value = 0.58142259850845401
$ws.Range("D10").Value = 0.54651145318371097
$ws.Range("D11").Value = 0.52963928012569705
$ws.Range("D12").Value = 0.50385732055235999
$ws.Range("D13").Value = 0.50368210152858095

# --- Column E ("Stopwords & Stemming") ---
$ws.Range("E3").Value = 0.619855248036006
$ws.Range("E4").Value = 0.61621888439964201
$ws.Range("E5").Value = 0.61232999551075395
$ws.Range("E6").Value = 0.57600736268756103
$ws.Range("E7").Value = 0.57034069602089399
$ws.Range("E8").Value = 0.56906147524167305
$ws.Range("E9").Value = 0.46975509395204801
$ws.Range("E10").Value = 0.44457320367430803
$ws.Range("E11").Value = 0.43007156267639701
$ws.Range("E12").Value = 0.40852024687831501
$ws.Range("E13").Value = 0.40806047676337198

# Make "Data" the active sheet with D15 selected (Chart1 was the active
# sheet before; now the data sheet is active with a fresh selection).
$ws.Activate()
$ws.Range("D15").Select()
